# Swap the contents of columns A and B for every used row, then restore
# the full-column selection over A:B (matches the author's "select A:B,
# swap columns" edit captured in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $colA = $ws.Cells.Item($r, 1)
    $colB = $ws.Cells.Item($r, 2)

    $valA = $colA.Value2
    $valB = $colB.Value2

    $colA.Value = $valB
    $colB.Value = $valA
}

$ws.Range("A1:B1048576").Select()
